# Update odds values on Sheet1 to reflect latest FlashScore data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.95
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 4.2
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5
$ws.Range("AC2").Value = 6.5
$ws.Range("BD2").Value = 126

# Row 4
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("Q4").Value = 2.1
$ws.Range("R4").Value = 1.7

# Row 5
$ws.Range("N5").Value = 8
$ws.Range("O5").Value = 1.4
$ws.Range("P5").Value = 3

# Row 6
$ws.Range("G6").Value = 1.75
$ws.Range("H6").Value = 3.2
$ws.Range("I6").Value = 5.25
$ws.Range("X6").Value = 7
$ws.Range("AF6").Value = 81
$ws.Range("AI6").Value = 23
$ws.Range("AK6").Value = 51
$ws.Range("AO6").Value = 10
$ws.Range("AW6").Value = 6.5
$ws.Range("AZ6").Value = 126
